$d = $word.ActiveDocument

# 1. Update the date paragraph (first paragraph of the document, outside the table).
#    Use Find/Replace scoped to that single paragraph's Range so nothing else is touched.
$dateRange = $d.Paragraphs(1).Range
$dateRange.Find.Execute("2025-11-05 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-11-06 Thursday", 2) | Out-Null

# 2. Update each multiplication-fact cell in the table.
#    Setting Cell.Range.Text directly (rather than Find.Execute on the cell Range)
#    keeps each edit strictly scoped to that one cell -- important because several
#    old/new fact strings collide across cells (e.g. "53×95=5035" is both an old value
#    in one cell and the new value written into another), and Find.Execute scoped to a
#    cell Range has been observed to still match/replace the first hit in the whole
#    story rather than staying within that cell.
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "18×98=1764"
$t.Cell(1, 2).Range.Text = "81×69=5589"
$t.Cell(1, 3).Range.Text = "87×45=3915"
$t.Cell(1, 4).Range.Text = "89×62=5518"
$t.Cell(1, 5).Range.Text = "26×12=312"
$t.Cell(5, 1).Range.Text = "12×48=576"
$t.Cell(5, 2).Range.Text = "42×43=1806"
$t.Cell(5, 3).Range.Text = "53×95=5035"
$t.Cell(5, 4).Range.Text = "93×70=6510"
$t.Cell(5, 5).Range.Text = "39×92=3588"
$t.Cell(10, 1).Range.Text = "72×52=3744"
$t.Cell(10, 2).Range.Text = "27×54=1458"
$t.Cell(10, 3).Range.Text = "42×33=1386"
$t.Cell(10, 4).Range.Text = "59×26=1534"
$t.Cell(10, 5).Range.Text = "88×81=7128"
$t.Cell(15, 1).Range.Text = "92×14=1288"
$t.Cell(15, 2).Range.Text = "14×91=1274"
$t.Cell(15, 3).Range.Text = "23×29=667"
$t.Cell(15, 4).Range.Text = "21×72=1512"
$t.Cell(15, 5).Range.Text = "25×56=1400"
$t.Cell(20, 1).Range.Text = "11×54=594"
$t.Cell(20, 2).Range.Text = "26×58=1508"
$t.Cell(20, 3).Range.Text = "22×72=1584"
$t.Cell(20, 4).Range.Text = "76×46=3496"
$t.Cell(20, 5).Range.Text = "84×97=8148"
